$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# --- Workbook window geometry (best effort; mirrors the saved window position/size) ---
$excel.ActiveWindow.Left = -120
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 19440
$excel.ActiveWindow.Height = 14880

# --- New row 12: DownloadsFolder config entry ---
$ws.Range("A12").Value = "DownloadsFolder"
$ws.Range("B12").Value = "C:\Users\RollLe01\Downloads\"
$ws.Range("C12").Value = "Download folder's path"

# --- Turn the RecipientTo / RecipientCC e-mail addresses (B10 / B11) into mailto: hyperlinks ---
# Insert the hyperlinks first (this is what produces the rId2 / rId3 relationships),
# then restore the existing "Hyperlink" cell style (copied from B3, which already uses it)
# so the cells keep the same visual style used elsewhere in the sheet.
$ws.Hyperlinks.Add($ws.Range("B11"), "mailto:lester.rollan@lexisnexisrisk.com")
$ws.Hyperlinks.Add($ws.Range("B10"), "mailto:lester.rollan@lexisnexisrisk.com")

$ws.Range("B3").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the selected cell shown when the workbook was last saved ---
$ws.Activate()
$ws.Range("B10").Select()
